# Export Groceris 2 - 06
# Appends a newly-exported grocery log batch (dated 2025-06-02, serial 45810)
# to both the "Log Per Recipe" and "Log Combined" worksheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Log Per Recipe")
$ws2 = $wb.Worksheets.Item("Log Combined")

# ---------------------------------------------------------------------------
# 1. "Log Per Recipe" - append rows 177-200
# ---------------------------------------------------------------------------
$recipeRows = @(
    @("CHLOETING-BROCCOLIRIJST", 11, "Tofu",         "TOFU",         1925,                "g"),
    @("CHLOETING-BROCCOLIRIJST", 11, "Broccoli",      "BROCCOLI",     2200,                "g"),
    @("CHLOETING-BROCCOLIRIJST", 11, "Lenteui",       "LENTEUI",      11,                  "u"),
    @("CHLOETING-BROCCOLIRIJST", 11, "Rijst",         "RIJST",        660,                 "g"),
    @("CHLOETING-BROCCOLIRIJST", 11, "Knoflookteen",  "KNOFLOOKTEEN", 11,                  "u"),
    @("CHLOETING-BROCCOLIRIJST", 11, "Soyasaus",      "SOYASAUS",     253,                 "g"),
    @("CHLOETING-BROCCOLIRIJST", 11, "Mirin",         "MIRIN",        242,                 "g"),
    @("PIZZAELS    ",             1, "Pizzadeeg",     "PIZZADEEG",    1,                   "u"),
    @("PIZZAELS    ",             1, "Passata",       "PASSATA",      75,                  "g"),
    @("PIZZAELS    ",             1, "Mozarella",     "MOZARELLA",    50,                  "g"),
    @("PIZZAJOE    ",             1, "Pizzadeeg",     "PIZZADEEG",    1,                   "u"),
    @("PIZZAJOE    ",             1, "Passata",       "PASSATA",      75,                  "g"),
    @("PIZZAJOE    ",             1, "Mozarella",     "MOZARELLA",    50,                  "g"),
    @("PIZZAJOE    ",             1, "Anjovis",       "ANJOVIS",      25,                  "g"),
    @("PIZZAJOE    ",             1, "Rode paprika",  "RODEPAPRIKA",  1,                   "u"),
    @("PIZZAJOE    ",             1, "Courgette",     "COURGETTE",    0.5,                 "u"),
    @("QUICHE      ",            12, "Bladerdeeg",    "BLADERDEEG",   2,                   "u"),
    @("QUICHE      ",            12, "Prei",          "PREI",         468,                 "g"),
    @("QUICHE      ",            12, "Volle Room",    "VOLLEROOM",    300,                 "g"),
    @("QUICHE      ",            12, "Mozarella",     "MOZARELLA",    180,                 "g"),
    @("QUICHE      ",            12, "Ei",            "EI",           9.6000000000000014,  "u"),
    @("QUICHE      ",            12, "Champignons",   "CHAMPIGNONS",  300,                 "g"),
    @("QUICHE      ",            12, "Miso",          "MISO",         180,                 "g"),
    @("QUICHE      ",            12, "Seitan",        "SEITAN",       600,                 "g")
)

$exportSerial = 45810
$firstRecipeRow = 177
$r = $firstRecipeRow
foreach ($row in $recipeRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $exportSerial
    $r = $r + 1
}
$lastRecipeRow = $r - 1

# Re-use the existing date-column style (built-in m/d/yyyy number format)
# instead of letting a brand new style get created.
$ws1.Range("G176").Copy()
$ws1.Range(("G" + $firstRecipeRow + ":G" + $lastRecipeRow)).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. "Log Combined" - append rows 154-187
# ---------------------------------------------------------------------------
$combinedRows = @(
    @("Ananas",          1,   "u"),
    @("Anjovis",          25,  "g"),
    @("Appel",            20,  "u"),
    @("Banaan",           2,   "u"),
    @("Bladerdeeg",       2,   "u"),
    @("Broccoli",         2200,"g"),
    @("Champignons",      300, "g"),
    @("Courgette",        0.5, "u"),
    @("Druiven",          2,   "u"),
    @("Edamame",          1200,"g"),
    @("Ei",               45.6,"u"),
    @("Knoflookteen",     11,  "u"),
    @("Lenteui",          11,  "u"),
    @("Melk",             3,   "u"),
    @("Mirin",            242, "g"),
    @("Miso",             180, "g"),
    @("Mozarella",        280, "g"),
    @("Passata",          150, "g"),
    @("Pindakaas",        1,   "u"),
    @("Pizzadeeg",        2,   "u"),
    @("Prei",             468, "g"),
    @("Rijst",            660, "g"),
    @("Rode paprika",     1,   "u"),
    @("Seitan",           600, "g"),
    @("Skyr",             4,   "u"),
    @("Soya melk",        4,   "u"),
    @("Soyasaus",         253, "g"),
    @("Tofu",             1925,"g"),
    @("Vegan Yoghurt",    2,   "u"),
    @("Volle Room",       300, "g"),
    @("Baguette",         2,   "u"),
    @("Rode Wijnazijn",   1,   "u"),
    @("Granaatappel",     1,   "u"),
    @("Peer",             2,   "u")
)

$firstCombinedRow = 154
$r = $firstCombinedRow
foreach ($row in $combinedRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $exportSerial
    $r = $r + 1
}
$lastCombinedRow = $r - 1

$ws2.Range("D153").Copy()
$ws2.Range(("D" + $firstCombinedRow + ":D" + $lastCombinedRow)).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. View state: "Log Combined" keeps its own scroll/selection but is no
#    longer the active tab; "Log Per Recipe" becomes the active/selected tab.
# ---------------------------------------------------------------------------
$ws2.Select()
$ws2.Range("E183").Select()

$ws1.Select()
$ws1.Range("G197").Select()
